# Update workbook to reflect a new sale of 238.35 for "MEZA FERNANDEZ JONATHAN ALEXIS"
# in the PIEDRA SINTERIZADA group for the month of "julio" (row 9 in both detail sheets),
# and propagate totals/summary values accordingly.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("L9").Value = 238.35
$wsGrupo.Range("L14").Value = "2 de 12"

# --- Sheet "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F9").Value = 238.35
$wsMensual.Range("F14").Value = 1252.07

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D2").Value = 5399.42
$wsCumpl.Range("E2").Value = -5399.42
$wsCumpl.Range("D4").Value = 10514.07
$wsCumpl.Range("E4").Value = 3209.27
$wsCumpl.Range("F4").Value = 0.766145122105843
